# feat: add 2022-Q3 data
#
# The "2022-Q2" sheet always holds the most-recent quarter's fund-holding
# data. When a new quarter arrives:
#   1. The current "2022-Q2" sheet is archived (copied) right after itself,
#      keeping its existing values, and the archived copy is renamed back
#      to "2022-Q2".
#   2. The original sheet is renamed to "2022-Q3" and its values are
#      overwritten with the new quarter's numbers.
#   3. The summary sheet ("总计") gets a new row inserted for 2022-Q3 (with
#      the same holding-count/market-value the prior quarter reported) and
#      every older row shifts down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fund-holding sheets: archive the current "2022-Q2" sheet, then turn
#    the original into the new "2022-Q3" sheet with updated figures.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Copy it to sit immediately after itself - this copy becomes the
# preserved historical "2022-Q2" snapshot.
$wsQ2.Copy($null, $wsQ2)
$wsQ2Archive = $wb.Worksheets.Item($wsQ2.Index + 1)
$wsQ2Archive.Name = "2022-Q2 archive"

# Rename the original sheet to the new quarter, then update its figures.
$wsQ2.Name = "2022-Q3"

$wsQ2.Range("D2").NumberFormat = "@"
$wsQ2.Range("D2").Value = "0.20"

$wsQ2.Range("E2").NumberFormat = "@"
$wsQ2.Range("E2").Value = "94.22"

$wsQ2.Range("F2").NumberFormat = "@"
$wsQ2.Range("F2").Value = "2.72"

$wsQ2.Range("G2").NumberFormat = "@"
$wsQ2.Range("G2").Value = "0.0054"

# Give the archived snapshot its proper name back now that the original
# no longer holds it.
$wsQ2Archive.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 2) Summary sheet ("总计"): insert a new row for 2022-Q3 above the old
#    2022-Q2 row, pushing every later quarter down by one row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

# The insert leaves column A on the new row without the shared style the
# rest of the column uses - copy it over from the row below first so the
# new row matches the rest of the table.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

# Write the final table contents explicitly (row 2 is the new 2022-Q3
# entry; rows 3-7 are the prior quarters, each shifted down one row).
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.01

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.16

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q4"
$wsTotal.Range("C5").Value = 3
$wsTotal.Range("D5").Value = 0.06

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q3"
$wsTotal.Range("C6").Value = 2
$wsTotal.Range("D6").Value = 0.08

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2020-Q4"
$wsTotal.Range("C7").Value = 2
$wsTotal.Range("D7").Value = 0.05
